# Included Fluid Menu List
# Adds a new "Fluid Menu List" navigation entry (with its own section
# header) to the Objects_Navigation worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objects_Navigation")

# --- New section header row (row 25), mirrors the existing section
# header rows (e.g. A22:G22) both in content pattern and formatting. ---
# Copy the formatting (bold white font on grey fill, centered) from an
# existing section-header row so the new header matches the others.
$ws.Range("A22:G22").Copy()
$ws.Range("A25:G25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Merge the header row across A:G like the other section headers.
$ws.Range("A25:G25").Merge()

$ws.Range("A25").Value = "User Site Fluid/Fulcrum Menu Navigation"

# --- New data row (row 26) describing the "Fluid Menu List" object. ---
$ws.Range("B26").Value = "Fluid Menu List"
$ws.Range("C26").Value = "xpath"
$ws.Range("D26").Value = "link"
$ws.Range("E26").Value = ".//span[contains(@class,'-item')]"

# --- Update the view selection to the newly added last row. ---
$ws.Activate()
$ws.Range("A26:XFD26").Select()
